# Generate Report for Handoff
# Updates the "Latest Handoff" timestamps (and Priority markers) for the
# six files that just had a new handoff package generated (rows 8-13 of
# each language report, corresponding to the files whose handoff xliff
# was just regenerated).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows 8 through 13 were refreshed by the new handoff generation.
$rows = 8..13

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-08-14 02:31:11"

    # zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority"
    $zhcn.Range("H$r").Value = "2016-08-14 02:31:00"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority"
    $dede.Range("H$r").Value = "2016-08-14 02:31:11"
    $dede.Range("E$r").Value = "ht"
}
